$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "Price" cells hold plain text that can look numeric (e.g. "1.003").
# Force text storage so values are not silently coerced to Double (which would
# drop formatting like trailing zeros, flip to scientific notation, etc.),
# then restore the original cell style so no visible formatting changes.

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.440.64'
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = '  -0.01%  '
$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.810.77'
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = '  +0.40%  '
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = '  -0.65%  '
$ws.Range("E5").Value = '  -0.43%  '
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '305.44'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  -0.80%  '
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4510'
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = '  -0.64%  '
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3581'
$ws.Range("D8").Style = $origStyle
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.28'
$ws.Range("D9").Style = $origStyle
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07055'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  -0.76%  '
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8885'
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  +1.61%  '
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07774'
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = '  +0.59%  '
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.30'
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = '  -0.03%  '
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.778.92'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  -3.92%  '
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.268'
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = '  +0.10%  '
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.297'
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = '  -0.62%  '
$ws.Range("E17").Value = '  -1.16%  '
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008513'
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = '  -0.65%  '
$ws.Range("E20").Value = '  -0.47%  '
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.485.50'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("E22").Value = '  -0.57%  '
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.947'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  -0.38%  '
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.033.03'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  -0.30%  '
$ws.Range("E25").Value = '  +0.89%  '
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.957'
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = '  -0.82%  '
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.13'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  +0.49%  '
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.76'
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = '  -0.76%  '
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.044'
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = '  +2.29%  '
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '112.14'
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = '  -0.28%  '
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.819'
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  +0.19%  '
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08682'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  +0.43%  '
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.135'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  +2.93%  '
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7451'
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = '  +2.48%  '
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.733'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  +6.97%  '
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.423'
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = '  +0.02%  '
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.108'
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = '  -0.27%  '
$ws.Range("E38").Value = '  -1.28%  '
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01926'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  +0.00%  '
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.899'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  +0.65%  '
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.05081'
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = '  -0.03%  '
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5076'
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = '  +1.67%  '
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.735'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  -2.81%  '
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1504'
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  -3.64%  '
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.034'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  -0.64%  '
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4711'
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = '  +2.64%  '
$ws.Range("E47").Value = '  -0.48%  '
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.01'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  +0.63%  '
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '100.05'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  -1.51%  '
$ws.Range("E50").Value = '  -0.80%  '
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05983'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  -0.03%  '
